# Replace/extend the data table in Sheet1 (A1:D7 -> A1:D17) with new
# values, per the target revision of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 0
$ws.Cells.Item(1, 2).Value = 109.1786205742166
$ws.Cells.Item(1, 3).Value = 10.37097476852108
$ws.Cells.Item(1, 4).Value = 1.233906563344711
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 112.511197023974
$ws.Cells.Item(2, 3).Value = 10.23179269695487
$ws.Cells.Item(2, 4).Value = 1.088284220444473
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 102.7756044855822
$ws.Cells.Item(3, 3).Value = 10.58245397251616
$ws.Cells.Item(3, 4).Value = 1.054717717983383
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 103.166746438103
$ws.Cells.Item(4, 3).Value = 10.8996482801407
$ws.Cells.Item(4, 4).Value = 0.7926528662174265
$ws.Cells.Item(5, 1).Value = 6
$ws.Cells.Item(5, 2).Value = 98.54296266253425
$ws.Cells.Item(5, 3).Value = 10.69053567380967
$ws.Cells.Item(5, 4).Value = 1.064762899394597
$ws.Cells.Item(6, 1).Value = 7
$ws.Cells.Item(6, 2).Value = 104.3718034604765
$ws.Cells.Item(6, 3).Value = 10.33323711615574
$ws.Cells.Item(6, 4).Value = 0.9435102513415738
$ws.Cells.Item(7, 1).Value = 9
$ws.Cells.Item(7, 2).Value = 101.1011844393747
$ws.Cells.Item(7, 3).Value = 10.37965696862676
$ws.Cells.Item(7, 4).Value = 1.066588561399274
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = 97.03337002844215
$ws.Cells.Item(8, 3).Value = 10.76193278594143
$ws.Cells.Item(8, 4).Value = 1.006258967614246
$ws.Cells.Item(9, 1).Value = 11
$ws.Cells.Item(9, 2).Value = 97.64265315771384
$ws.Cells.Item(9, 3).Value = 10.29382548146205
$ws.Cells.Item(9, 4).Value = 1.124525108184251
$ws.Cells.Item(10, 1).Value = 12
$ws.Cells.Item(10, 2).Value = 99.46655366285279
$ws.Cells.Item(10, 3).Value = 10.52890109413067
$ws.Cells.Item(10, 4).Value = 0.8283571727685592
$ws.Cells.Item(11, 1).Value = 14
$ws.Cells.Item(11, 2).Value = 93.15156922950189
$ws.Cells.Item(11, 3).Value = 10.58962240117504
$ws.Cells.Item(11, 4).Value = 0.858600470404917
$ws.Cells.Item(12, 1).Value = 15
$ws.Cells.Item(12, 2).Value = 92.78555524585381
$ws.Cells.Item(12, 3).Value = 10.13239105666403
$ws.Cells.Item(12, 4).Value = 1.09887750072958
$ws.Cells.Item(13, 1).Value = 17
$ws.Cells.Item(13, 2).Value = 88.81350597853542
$ws.Cells.Item(13, 3).Value = 10.06948390580832
$ws.Cells.Item(13, 4).Value = 1.109296649182347
$ws.Cells.Item(14, 1).Value = 19
$ws.Cells.Item(14, 2).Value = 90.42192385422325
$ws.Cells.Item(14, 3).Value = 10.15493461409919
$ws.Cells.Item(14, 4).Value = 0.9148563577195951
$ws.Cells.Item(15, 1).Value = 20
$ws.Cells.Item(15, 2).Value = 89.38019988050237
$ws.Cells.Item(15, 3).Value = 10.45565062792593
$ws.Cells.Item(15, 4).Value = 0.755169316139251
$ws.Cells.Item(16, 1).Value = 24
$ws.Cells.Item(16, 2).Value = 84.6488430913671
$ws.Cells.Item(16, 3).Value = 10.11986427361115
$ws.Cells.Item(16, 4).Value = 0.9047655866026827
$ws.Cells.Item(17, 1).Value = 27
$ws.Cells.Item(17, 2).Value = 81.35088578706497
$ws.Cells.Item(17, 3).Value = 10.21020303845059
$ws.Cells.Item(17, 4).Value = 0.9203340417559835
